$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.024" or
# "28.409.18" are not auto-converted to numbers by Excel, matching the
# original inline-string representation of these cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.409.18"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "1.867.96"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("D4").Value = "1.024"
$ws.Range("E4").Value = "  +2.25%  "

$ws.Range("D5").Value = "317.08"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").Value = "1.021"
$ws.Range("E6").Value = "  +2.15%  "

$ws.Range("D7").Value = "0.5110"
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").Value = "0.3961"
$ws.Range("E8").Value = "  +2.28%  "

$ws.Range("D9").Value = "0.08338"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "42.04"
$ws.Range("E10").Value = "  +1.48%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "1.108"
$ws.Range("E11").Value = "  -0.56%  "

$ws.Range("D12").Value = "6.240"
$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("D13").Value = "20.42"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("D14").Value = "1.833.11"
$ws.Range("E14").Value = "  -1.33%  "

$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "1.024"
$ws.Range("E15").Value = "  +2.40%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "7.211"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").Value = "0.00001107"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").Value = "90.90"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").Value = "0.06767"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.021"
$ws.Range("E20").Value = "  +2.16%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "17.66"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").Value = "5.945"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("D23").Value = "28.437.49"
$ws.Range("E23").Value = "  +1.27%  "

$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  +0.42%  "

$ws.Range("D25").Value = "2.291"
$ws.Range("E25").Value = "  +2.72%  "

$ws.Range("D26").Value = "161.61"
$ws.Range("E26").Value = "  +2.14%  "

$ws.Range("D27").Value = "2.030.15"
$ws.Range("E27").Value = "  -2.04%  "

$ws.Range("D28").Value = "20.78"
$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("D29").Value = "2.356"
$ws.Range("E29").Value = "  -4.23%  "

$ws.Range("D30").Value = "127.30"
$ws.Range("E30").Value = "  +2.13%  "

$ws.Range("D31").Value = "0.1048"
$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("D32").Value = "1.030"
$ws.Range("E32").Value = "  -0.28%  "

$ws.Range("D33").Value = "5.772"
$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("D34").Value = "3.635"
$ws.Range("E34").Value = "  +1.41%  "

$ws.Range("D35").Value = "0.02421"
$ws.Range("E35").Value = "  -0.54%  "

$ws.Range("D36").Value = "0.06470"
$ws.Range("E36").Value = "  -0.88%  "

$ws.Range("D37").Value = "0.2181"
$ws.Range("E37").Value = "  -0.78%  "

$ws.Range("D38").Value = "8.885"
$ws.Range("E38").Value = "  -7.04%  "

$ws.Range("D39").Value = "1.271"
$ws.Range("E39").Value = "  +3.61%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.180"
$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6398"
$ws.Range("E41").Value = "  -1.08%  "

$ws.Range("D42").Value = "5.009"
$ws.Range("E42").Value = "  +1.28%  "

$ws.Range("D43").Value = "11.19"
$ws.Range("E43").Value = "  -0.22%  "

$ws.Range("D44").Value = "0.6008"
$ws.Range("E44").Value = "  -1.29%  "

$ws.Range("D45").Value = "13.06"
$ws.Range("E45").Value = "  +0.22%  "

$ws.Range("D46").Value = "3.722"
$ws.Range("E46").Value = "  +1.98%  "

$ws.Range("D47").Value = "1.218"
$ws.Range("E47").Value = "  -4.21%  "

$ws.Range("D48").Value = "1.986"
$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("D49").Value = "121.83"
$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("D50").Value = "1.204"
$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("D51").Value = "0.06854"
$ws.Range("E51").Value = "  -0.50%  "

# Restore the original (default/"Normal") cell style on column D now that
# the text values are safely stored, so no stray number-format styling is
# left behind on the cells.
$priceRange.Style = "Normal"
